$wb = $excel.ActiveWorkbook

# --- Rushing sheet updates ---
$rushing = $wb.Worksheets.Item("Rushing")
$rushing.Range("D2").Value = 13
$rushing.Range("E2").Value = 13

$rushing.Range("C4").Value = 202
$rushing.Range("D4").Value = 99
$rushing.Range("F4").Value = 42

$rushing.Range("C6").Value = 7

$rushing.Range("C9").Value = 6

# --- Receiving sheet updates ---
$receiving = $wb.Worksheets.Item("Receiving")

$receiving.Range("C2").Value = 55
$receiving.Range("D2").Value = 47
$receiving.Range("E2").Value = 5
$receiving.Range("F2").Value = 5

$receiving.Range("C3").Value = 32
$receiving.Range("D3").Value = 27

$receiving.Range("C5").Value = 84
$receiving.Range("D5").Value = 63

$receiving.Range("C6").Value = 100
$receiving.Range("D6").Value = 72
$receiving.Range("E6").Value = 34
$receiving.Range("F6").Value = 21

$receiving.Range("C7").Value = 98
$receiving.Range("D7").Value = 70
$receiving.Range("E7").Value = 45
$receiving.Range("F7").Value = 24

$receiving.Range("C11").Value = 67
$receiving.Range("D11").Value = 53
$receiving.Range("E11").Value = 10
$receiving.Range("F11").Value = 9
$receiving.Range("G11").Value = 5
$receiving.Range("H11").Value = 3
